$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text so values like "1.00" or "6.20"
# are not silently reinterpreted as numbers, matching the source data which
# stores prices as literal strings (inline strings in the original file).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '64.299.69'
$ws.Range("E2").Value = '  +2.28%  '

# Row 3
$ws.Range("D3").Value = '3.076.13'
$ws.Range("E3").Value = '  +1.08%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").Value = '559.09'
$ws.Range("E5").Value = '  +1.97%  '

# Row 6
$ws.Range("D6").Value = '145.21'
$ws.Range("E6").Value = '  +5.34%  '

# Row 7
$ws.Range("E7").Value = '  +0.12%  '

# Row 8
$ws.Range("D8").Value = '3.072.88'
$ws.Range("E8").Value = '  +1.09%  '

# Row 9
$ws.Range("E9").Value = '  +1.04%  '

# Row 10
$ws.Range("E10").Value = '  +2.55%  '

# Row 11
$ws.Range("D11").Value = '6.20'
$ws.Range("E11").Value = '  -1.95%  '

# Row 12
$ws.Range("D12").Value = '0.469'
$ws.Range("E12").Value = '  +4.39%  '

# Row 13
$ws.Range("E13").Value = '  +0.98%  '

# Row 14
$ws.Range("D14").Value = '35.14'
$ws.Range("E14").Value = '  +2.11%  '

# Row 15
$ws.Range("D15").Value = '3.578.47'
$ws.Range("E15").Value = '  +1.28%  '

# Row 16
$ws.Range("D16").Value = '64.371.59'
$ws.Range("E16").Value = '  +2.35%  '

# Row 17
$ws.Range("D17").Value = '3.076.34'
$ws.Range("E17").Value = '  +1.25%  '

# Row 18
$ws.Range("E18").Value = '  +1.60%  '

# Row 19
$ws.Range("D19").Value = '6.75'
$ws.Range("E19").Value = '  +0.86%  '

# Row 20
$ws.Range("D20").Value = '477.80'
$ws.Range("E20").Value = '  -0.21%  '

# Row 21
$ws.Range("D21").Value = '13.87'
$ws.Range("E21").Value = '  +2.21%  '

# Row 22
$ws.Range("E22").Value = '  +1.01%  '

# Row 23
$ws.Range("D23").Value = '7.54'
$ws.Range("E23").Value = '  +5.40%  '

# Row 24
$ws.Range("D24").Value = '13.60'
$ws.Range("E24").Value = '  +10.06%  '

# Row 25
$ws.Range("D25").Value = '81.30'
$ws.Range("E25").Value = '  +1.26%  '

# Row 26
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.07%  '

# Row 27
$ws.Range("E27").Value = '  +2.38%  '

# Row 28
$ws.Range("D28").Value = '8.02'
$ws.Range("E28").Value = '  +2.36%  '

# Row 29
$ws.Range("E29").Value = '  +4.96%  '

# Row 30
$ws.Range("E30").Value = '  -0.03%  '

# Row 31
$ws.Range("D31").Value = '26.12'
$ws.Range("E31").Value = '  +1.24%  '

# Row 32
$ws.Range("D32").Value = '1.15'
$ws.Range("E32").Value = '  +0.63%  '

# Row 33
$ws.Range("D33").Value = '2.49'
$ws.Range("E33").Value = '  +4.25%  '

# Row 34
$ws.Range("D34").Value = '5.57'
$ws.Range("E34").Value = '  -1.26%  '

# Row 35
$ws.Range("D35").Value = '55.83'
$ws.Range("E35").Value = '  +1.08%  '

# Row 36
$ws.Range("D36").Value = '6.17'
$ws.Range("E36").Value = '  +4.33%  '

# Row 37
$ws.Range("D37").Value = '457.14'
$ws.Range("E37").Value = '  -0.72%  '

# Row 38
$ws.Range("D38").Value = '3.02'
$ws.Range("E38").Value = '  +19.29%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.0405'
$ws.Range("E39").Value = '  +3.41%  '

# Row 40
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.0824'
$ws.Range("E40").Value = '  +1.87%  '

# Row 41
$ws.Range("D41").Value = '2.958.03'
$ws.Range("E41").Value = '  -3.34%  '

# Row 42
$ws.Range("D42").Value = '8.23'
$ws.Range("E42").Value = '  +0.33%  '

# Row 43
$ws.Range("E43").Value = '  -1.56%  '

# Row 44
$ws.Range("D44").Value = '27.81'
$ws.Range("E44").Value = '  -0.85%  '

# Row 45
$ws.Range("D45").Value = '0.261'
$ws.Range("E45").Value = '  +4.15%  '

# Row 46
$ws.Range("D46").Value = '2.14'
$ws.Range("E46").Value = '  +5.39%  '

# Row 47
$ws.Range("E47").Value = '  +0.04%  '

# Row 48
$ws.Range("E48").Value = '  +2.73%  '

# Row 49
$ws.Range("D49").Value = '121.29'
$ws.Range("E49").Value = '  +4.38%  '

# Row 50
$ws.Range("D50").Value = '0.0₃0515'
$ws.Range("E50").Value = '  +2.42%  '

# Row 51
$ws.Range("D51").Value = '2.07'
$ws.Range("E51").Value = '  +1.39%  '
